# Update "Critical issues" sheet (GSC export) with refreshed coverage data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Critical issues")

# Header row (unchanged, kept for clarity)
$ws.Range("A1").Value = "Reason"
$ws.Range("B1").Value = "Source"
$ws.Range("C1").Value = "Validation"
$ws.Range("D1").Value = "Pages"

# Data rows 2-10, reordered/updated per latest GSC coverage export
$ws.Range("A2").Value = "Excluded by ‘noindex’ tag"
$ws.Range("B2").Value = "Website"
$ws.Range("C2").Value = "Not Started"
$ws.Range("D2").Value = 17.0

$ws.Range("A3").Value = "Blocked by robots.txt"
$ws.Range("B3").Value = "Website"
$ws.Range("C3").Value = "Not Started"
$ws.Range("D3").Value = 1.0

$ws.Range("A4").Value = "Alternate page with proper canonical tag"
$ws.Range("B4").Value = "Website"
$ws.Range("C4").Value = "Started"
$ws.Range("D4").Value = 62.0

$ws.Range("A5").Value = "Not found (404)"
$ws.Range("B5").Value = "Website"
$ws.Range("C5").Value = "Started"
$ws.Range("D5").Value = 51.0

$ws.Range("A6").Value = "Page with redirect"
$ws.Range("B6").Value = "Website"
$ws.Range("C6").Value = "Started"
$ws.Range("D6").Value = 15.0

$ws.Range("A7").Value = "Duplicate, Google chose different canonical than user"
$ws.Range("B7").Value = "Google systems"
$ws.Range("C7").Value = "Started"
$ws.Range("D7").Value = 25.0

$ws.Range("A8").Value = "Crawled - currently not indexed"
$ws.Range("B8").Value = "Google systems"
$ws.Range("C8").Value = "Started"
$ws.Range("D8").Value = 9.0

$ws.Range("A9").Value = "Server error (5xx)"
$ws.Range("B9").Value = "Website"
$ws.Range("C9").Value = "Passed"
$ws.Range("D9").Value = 0.0

$ws.Range("A10").Value = "Discovered - currently not indexed"
$ws.Range("B10").Value = "Google systems"
$ws.Range("C10").Value = "Passed"
$ws.Range("D10").Value = 0.0
